# --- "sus final na tese" --------------------------------------------------
# Renames the original sheet to "SUS Prototipo", duplicates it into a new
# sheet "SUS HMI" (becoming the active/selected tab) with updated SUS scores
# and an extra raw-data block (columns I:S) highlighting the per-user totals
# with a "Correto"/Good conditional style.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1. Rename the existing sheet and duplicate it to build the new sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SUS Prototipo"
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "SUS HMI"

# 2. Update the scores that differ on the new "SUS HMI" sheet.
$ws2.Range("G2").Value = 4
$ws2.Range("G3").Value = 2
$ws2.Range("E4").Value = 5
$ws2.Range("D5").Value = 2
$ws2.Range("E5").Value = 1
$ws2.Range("G5").Value = 1
$ws2.Range("E6").Value = 4
$ws2.Range("F7").Value = 2
$ws2.Range("E10").Value = 4
$ws2.Range("F10").Value = 5
$ws2.Range("D12").Value = 92.5
$ws2.Range("E12").Value = 87.5
$ws2.Range("G12").Value = 87.5
$ws2.Range("E14").Value = 89.375

# 3. Hide the raw-total helper column (C), as on the prototype sheet, but
#    hidden here.
$ws2.Columns.Item(3).Hidden = $true

# 4. Add the per-user raw SUS answers block in columns I:S.
$ws2.Range("I2").Value = "U2"
$ws2.Range("I3").Value = "U1"
$ws2.Range("I4").Value = "U3"
$ws2.Range("I5").Value = "U4"

$rawData = @{
    2 = @(4,1,5,1,4,1,4,2,4,1)
    3 = @(4,1,5,2,5,1,5,1,5,2)
    4 = @(4,1,5,2,4,2,5,1,5,1)
    5 = @(4,2,4,1,4,1,4,1,5,1)
}

foreach ($row in $rawData.Keys) {
    $values = $rawData[$row]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $col = 10 + $i   # column J = 10
        $ws2.Cells.Item($row, $col).Value = $values[$i]
    }
}

# 5. Apply the green "Correto" highlight style to the J2:S3 block.
$ws2.Range("J2:S3").Style = "Good"

# 6. Fix up selections / active sheet so that "SUS HMI" ends up the
#    tab that is selected and active, matching the authored workbook.
$ws1.Activate()
$ws1.Range("C1").Select()

$ws2.Activate()
$ws2.Range("N31").Select()
